# Updated symbol list on Tue Dec 27 20:36:55 UTC 2022 with GitHub Actions
# Applies the refreshed "Price" (column D) values and the two
# "Volume(1h)" (column E) label tweaks from the upstream scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    # Force the cell to Text format before/while assigning so a
    # numeric-looking string (e.g. "0.001554") is kept verbatim instead
    # of being parsed into a Double (which would drop trailing zeros /
    # switch to scientific notation). ClearFormats afterwards removes the
    # temporary "@" number format again so no stray style id is left on
    # the cell, matching the original (unstyled) cells.
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

# Column D ("Price") updates
Set-TextValue "D2"  "245.91"
Set-TextValue "D3"  "23.98"
Set-TextValue "D4"  "5.369"
Set-TextValue "D5"  "0.05832"
Set-TextValue "D6"  "6.470"
Set-TextValue "D7"  "3.345"
Set-TextValue "D8"  "0.8091"
Set-TextValue "D9"  "0.9202"
Set-TextValue "D11" "0.07374"
Set-TextValue "D12" "0.03108"
Set-TextValue "D14" "0.09369"
Set-TextValue "D15" "3.864"
Set-TextValue "D16" "0.001554"
Set-TextValue "D17" "0.04693"
Set-TextValue "D19" "0.006195"
Set-TextValue "D21" "0.004692"
Set-TextValue "D22" "0.00008798"
Set-TextValue "D28" "0.0002348"
Set-TextValue "D40" "0.03847"
Set-TextValue "D41" "0.003094"
Set-TextValue "D42" "0.1067"
Set-TextValue "D43" "0.002749"
Set-TextValue "D44" "0.008524"
Set-TextValue "D45" "0.00005250"
Set-TextValue "D47" "0.6856"

# Column E ("Volume(1h)") label updates
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"
